$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6520
$ws.Range("L3").Value = 7033
$ws.Range("B4").Value = 1721
$ws.Range("L4").Value = 1754
$ws.Range("L5").Value = 416
$ws.Range("L6").Value = 5766
$ws.Range("B7").Value = 23353
$ws.Range("L7").Value = 21489

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L2").Value = 73
$ws.Range("L7").Value = 242

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 434
$ws.Range("L3").Value = 498
$ws.Range("L6").Value = 342
$ws.Range("L7").Value = 1417

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 151
$ws.Range("L3").Value = 188
$ws.Range("L7").Value = 473

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 341
$ws.Range("L6").Value = 271
$ws.Range("L7").Value = 965

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L5").Value = 24
$ws.Range("L6").Value = 213
$ws.Range("L7").Value = 828

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 162
$ws.Range("L7").Value = 428

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L6").Value = 78
$ws.Range("L7").Value = 370

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 194
$ws.Range("L4").Value = 78
$ws.Range("L7").Value = 689
$ws.Range("L8").Value = 1417
$ws.Range("L15").Value = 186
$ws.Range("L19").Value = 594
$ws.Range("L20").Value = 540
$ws.Range("L31").Value = 215
$ws.Range("L32").Value = 26
$ws.Range("L33").Value = 965
$ws.Range("L37").Value = 828
$ws.Range("L43").Value = 163
$ws.Range("L47").Value = 149
$ws.Range("L48").Value = 278
$ws.Range("L51").Value = 265
$ws.Range("L52").Value = 455
$ws.Range("L53").Value = 242
$ws.Range("B63").Value = 425
$ws.Range("L63").Value = 66
$ws.Range("L65").Value = 428
$ws.Range("L76").Value = 340
$ws.Range("L78").Value = 282
$ws.Range("L79").Value = 599
$ws.Range("L83").Value = 473
$ws.Range("L85").Value = 1071
$ws.Range("L90").Value = 231
$ws.Range("L93").Value = 110
$ws.Range("L96").Value = 236
$ws.Range("L99").Value = 370
$ws.Range("B101").Value = 23353
$ws.Range("L101").Value = 21489

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 215

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L4").Value = 53
$ws.Range("L7").Value = 278

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 212
$ws.Range("L3").Value = 180
$ws.Range("L5").Value = 9
$ws.Range("L6").Value = 164
$ws.Range("L7").Value = 594

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 67
$ws.Range("L7").Value = 340

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L6").Value = 81
$ws.Range("L7").Value = 282

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 67
$ws.Range("L7").Value = 236

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 195
$ws.Range("L6").Value = 160
$ws.Range("L7").Value = 599

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 189
$ws.Range("L7").Value = 540

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L3").Value = 31
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 240
$ws.Range("L3").Value = 221
$ws.Range("L7").Value = 689

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 56
$ws.Range("L7").Value = 149

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L3").Value = 59
$ws.Range("L6").Value = 38
$ws.Range("L7").Value = 186

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L2").Value = 65
$ws.Range("L7").Value = 194

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 79
$ws.Range("L3").Value = 66
$ws.Range("L7").Value = 231

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L4").Value = 40
$ws.Range("L7").Value = 265

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L3").Value = 54
$ws.Range("L4").Value = 27
$ws.Range("L7").Value = 163

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 444
$ws.Range("L7").Value = 1071

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L6").Value = 127
$ws.Range("L7").Value = 455

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 78
